$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'FAPs'
$ws.Range('B2').Value = 'Ntf5'
$ws.Range('C2').Value = 'Ntrk1'
$ws.Range('D2').Value = 'ECs'
$ws.Range('E2').Value = 2
$ws.Range('F2').Value = 0.6666666666666666
$ws.Range('G2').Value = 0.3590443333333333
$ws.Range('H2').Value = 1.077133
$ws.Range('I2').Value = 0.480597973884934
$ws.Range('J2').Value = 0.480597973884934
$ws.Range('K2').Value = 1
$ws.Range('L2').Value = 0.3333333333333333
$ws.Range('M2').Value = 0.05351466666666666
$ws.Range('N2').Value = 0.160544
$ws.Range('O2').Value = 0.4186447970585551
$ws.Range('P2').Value = 0.4186447970585551
$ws.Range('Q2').Value = 0.01921413781688889
$ws.Range('R2').Value = 0.172927240352
$ws.Range('S2').Value = 0.201199841243811
$ws.Range('T2').Value = 0.201199841243811

# Row 3
$ws.Range('A3').Value = 'FAPs'
$ws.Range('B3').Value = 'Ntf5'
$ws.Range('C3').Value = 'Ntrk1'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('E3').Value = 2
$ws.Range('F3').Value = 0.6666666666666666
$ws.Range('G3').Value = 0.3590443333333333
$ws.Range('H3').Value = 1.077133
$ws.Range('I3').Value = 0.480597973884934
$ws.Range('J3').Value = 0.480597973884934
$ws.Range('K3').Value = 2
$ws.Range('L3').Value = 0.6666666666666666
$ws.Range('M3').Value = 0.05507833333333334
$ws.Range('N3').Value = 0.165235
$ws.Range('O3').Value = 0.4308773485273219
$ws.Range('P3').Value = 0.4308773485273218
$ws.Range('Q3').Value = 0.01977556347277778
$ws.Range('R3').Value = 0.177980071255
$ws.Range('S3').Value = 0.2070787806951434
$ws.Range('T3').Value = 0.2070787806951434

# Row 4
$ws.Range('A4').Value = 'FAPs'
$ws.Range('B4').Value = 'Ntf5'
$ws.Range('C4').Value = 'Ntrk1'
$ws.Range('D4').Value = 'MuSCs'
$ws.Range('E4').Value = 2
$ws.Range('F4').Value = 0.6666666666666666
$ws.Range('G4').Value = 0.3590443333333333
$ws.Range('H4').Value = 1.077133
$ws.Range('I4').Value = 0.480597973884934
$ws.Range('J4').Value = 0.480597973884934
$ws.Range('K4').Value = 1
$ws.Range('L4').Value = 0.3333333333333333
$ws.Range('M4').Value = 0.01923533333333333
$ws.Range('N4').Value = 0.057706
$ws.Range('O4').Value = 0.1504778544141231
$ws.Range('P4').Value = 0.1504778544141231
$ws.Range('Q4').Value = 0.006906337433111111
$ws.Range('R4').Value = 0.06215703689799999
$ws.Range('S4').Value = 0.07231935194597963
$ws.Range('T4').Value = 0.07231935194597963

# Row 5
$ws.Range('A5').Value = 'Inflammatory-Mac'
$ws.Range('B5').Value = 'Ntf5'
$ws.Range('C5').Value = 'Ntrk1'
$ws.Range('D5').Value = 'ECs'
$ws.Range('E5').Value = 1
$ws.Range('F5').Value = 0.3333333333333333
$ws.Range('G5').Value = 0.388034
$ws.Range('H5').Value = 1.164102
$ws.Range('I5').Value = 0.519402026115066
$ws.Range('J5').Value = 0.5194020261150661
$ws.Range('K5').Value = 1
$ws.Range('L5').Value = 0.3333333333333333
$ws.Range('M5').Value = 0.05351466666666666
$ws.Range('N5').Value = 0.160544
$ws.Range('O5').Value = 0.4186447970585551
$ws.Range('P5').Value = 0.4186447970585551
$ws.Range('Q5').Value = 0.02076551016533333
$ws.Range('R5').Value = 0.186889591488
$ws.Range('S5').Value = 0.2174449558147442
$ws.Range('T5').Value = 0.2174449558147442

# Row 6
$ws.Range('A6').Value = 'Inflammatory-Mac'
$ws.Range('B6').Value = 'Ntf5'
$ws.Range('C6').Value = 'Ntrk1'
$ws.Range('D6').Value = 'FAPs'
$ws.Range('E6').Value = 1
$ws.Range('F6').Value = 0.3333333333333333
$ws.Range('G6').Value = 0.388034
$ws.Range('H6').Value = 1.164102
$ws.Range('I6').Value = 0.519402026115066
$ws.Range('J6').Value = 0.5194020261150661
$ws.Range('K6').Value = 2
$ws.Range('L6').Value = 0.6666666666666666
$ws.Range('M6').Value = 0.05507833333333334
$ws.Range('N6').Value = 0.165235
$ws.Range('O6').Value = 0.4308773485273219
$ws.Range('P6').Value = 0.4308773485273218
$ws.Range('Q6').Value = 0.02137226599666667
$ws.Range('R6').Value = 0.19235039397
$ws.Range('S6').Value = 0.2237985678321785
$ws.Range('T6').Value = 0.2237985678321785

# Row 7
$ws.Range('A7').Value = 'Inflammatory-Mac'
$ws.Range('B7').Value = 'Ntf5'
$ws.Range('C7').Value = 'Ntrk1'
$ws.Range('D7').Value = 'MuSCs'
$ws.Range('E7').Value = 1
$ws.Range('F7').Value = 0.3333333333333333
$ws.Range('G7').Value = 0.388034
$ws.Range('H7').Value = 1.164102
$ws.Range('I7').Value = 0.519402026115066
$ws.Range('J7').Value = 0.5194020261150661
$ws.Range('K7').Value = 1
$ws.Range('L7').Value = 0.3333333333333333
$ws.Range('M7').Value = 0.01923533333333333
$ws.Range('N7').Value = 0.057706
$ws.Range('O7').Value = 0.1504778544141231
$ws.Range('P7').Value = 0.1504778544141231
$ws.Range('Q7').Value = 0.007463963334666666
$ws.Range('R7').Value = 0.067175670012
$ws.Range('S7').Value = 0.07815850246814349
$ws.Range('T7').Value = 0.07815850246814349
